$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2654.1843
$ws.Range("I15").Value = 2654.1843
$ws.Range("K15").Value = 7962.5529
$ws.Range("M15").Value = -7793.5529

# Row 32
$ws.Range("H32").Value = 125002160
$ws.Range("I32").Value = 333334270
$ws.Range("J32").Value = 2900
$ws.Range("K32").Value = 333334270
$ws.Range("L32").Value = 2900
$ws.Range("M32").Value = -333333944
$ws.Range("N32").Value = -3552

# Row 86
$ws.Range("H86").Value = 4152.923
$ws.Range("I86").Value = 3524
$ws.Range("J86").Value = 4432.4443
$ws.Range("K86").Value = 3524
$ws.Range("L86").Value = 4432.4443
$ws.Range("M86").Value = -2401
$ws.Range("N86").Value = -6678.4443

# Row 89
$ws.Range("H89").Value = 4152.923
$ws.Range("I89").Value = 3524
$ws.Range("J89").Value = 4432.4443
$ws.Range("K89").Value = 17620
$ws.Range("L89").Value = 22162.2215
$ws.Range("M89").Value = -12004
$ws.Range("N89").Value = -33394.2215

# Row 92
$ws.Range("H92").Value = 1528.1177
$ws.Range("I92").Value = 1677
$ws.Range("J92").Value = 833.3333
$ws.Range("K92").Value = 1677
$ws.Range("L92").Value = 833.3333
$ws.Range("M92").Value = -429
$ws.Range("N92").Value = -3329.3333

# Row 113
$ws.Range("H113").Value = 2807
$ws.Range("I113").Value = 3550.625
$ws.Range("J113").Value = 1957.1428
$ws.Range("K113").Value = 3550.625
$ws.Range("L113").Value = 1957.1428
$ws.Range("M113").Value = -296.625
$ws.Range("N113").Value = -8465.1428

# Row 121
$ws.Range("H121").Value = 608.17145
$ws.Range("J121").Value = 596.5454999999999
$ws.Range("L121").Value = 1789.6365
$ws.Range("N121").Value = -5283.6365

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 16995.518
$ws.Range("I32").Value = 14031.071
$ws.Range("K32").Value = 14031.071
$ws.Range("M32").Value = -13744.071

# Row 34
$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 10000
$ws.Range("K34").Value = 10000
$ws.Range("M34").Value = -9729

# Row 45
$ws.Range("H45").Value = 8473.058999999999
$ws.Range("I45").Value = 11612.8
$ws.Range("J45").Value = 3987.7144
$ws.Range("K45").Value = 11612.8
$ws.Range("L45").Value = 3987.7144
$ws.Range("M45").Value = -11235.8
$ws.Range("N45").Value = -4741.7144

# Row 74
$ws.Range("H74").Value = 817.6667
$ws.Range("I74").Value = 837.6
$ws.Range("J74").Value = 792.75
$ws.Range("K74").Value = 837.6
$ws.Range("L74").Value = 792.75
$ws.Range("M74").Value = 36.39999999999998
$ws.Range("N74").Value = -2540.75

# Row 77
$ws.Range("H77").Value = 817.6667
$ws.Range("I77").Value = 837.6
$ws.Range("J77").Value = 792.75
$ws.Range("K77").Value = 4188
$ws.Range("L77").Value = 3963.75
$ws.Range("M77").Value = 180
$ws.Range("N77").Value = -12699.75

# Row 97
$ws.Range("H97").Value = 1098.8889
$ws.Range("I97").Value = 986.25
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 986.25
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -490.25
$ws.Range("N97").Value = -2992

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3091.5
$ws.Range("I105").Value = 3173.5715
$ws.Range("J105").Value = 2900
$ws.Range("K105").Value = 3173.5715
$ws.Range("L105").Value = 2900
$ws.Range("M105").Value = -1426.5715
$ws.Range("N105").Value = -6394

# Row 107
$ws.Range("H107").Value = 1099.1052
$ws.Range("I107").Value = 1163.4615
$ws.Range("J107").Value = 959.6667
$ws.Range("K107").Value = 1163.4615
$ws.Range("L107").Value = 959.6667
$ws.Range("M107").Value = 756.5385000000001
$ws.Range("N107").Value = -4799.6667

# Row 126
$ws.Range("H126").Value = 51320
$ws.Range("J126").Value = 51320
$ws.Range("L126").Value = 51320
$ws.Range("N126").Value = -61200

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 4780
$ws.Range("I99").Value = 4000
$ws.Range("J99").Value = 4975
$ws.Range("K99").Value = 4000
$ws.Range("L99").Value = 4975
$ws.Range("M99").Value = -2502
$ws.Range("N99").Value = -7971

# Row 122
$ws.Range("H122").Value = 5000506
$ws.Range("I122").Value = 10000012
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 30000036
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -29997586
$ws.Range("N122").Value = -7900

# Row 126
$ws.Range("H126").Value = 4780
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 4975
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 14925
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -19865

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 202980
$ws.Range("I55").Value = 1000000
$ws.Range("J55").Value = 3725
$ws.Range("K55").Value = 3000000
$ws.Range("L55").Value = 11175
$ws.Range("M55").Value = -2999823
$ws.Range("N55").Value = -11529

# Row 68
$ws.Range("H68").Value = 982.48486
$ws.Range("I68").Value = 577.4091
$ws.Range("J68").Value = 1306.5454
$ws.Range("K68").Value = 1732.2273
$ws.Range("L68").Value = 3919.6362
$ws.Range("M68").Value = -921.2273
$ws.Range("N68").Value = -5541.6362

# Row 71
$ws.Range("H71").Value = 982.48486
$ws.Range("I71").Value = 577.4091
$ws.Range("J71").Value = 1306.5454
$ws.Range("K71").Value = 5196.6819
$ws.Range("L71").Value = 11758.9086
$ws.Range("M71").Value = -1140.6819
$ws.Range("N71").Value = -19870.9086

# Row 113
$ws.Range("H113").Value = 423.4
$ws.Range("I113").Value = 373.64285
$ws.Range("J113").Value = 431.5
$ws.Range("K113").Value = 1120.92855
$ws.Range("L113").Value = 1294.5
$ws.Range("M113").Value = 1049.07145
$ws.Range("N113").Value = -5634.5

# Row 132
$ws.Range("H132").Value = 632.35297
$ws.Range("I132").Value = 596.6667
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 5370.0003
$ws.Range("L132").Value = 8100
$ws.Range("M132").Value = -2840.0003
$ws.Range("N132").Value = -13160

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3946.5386
$ws.Range("I80").Value = 4025.4167
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 4025.4167
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -3027.4167
$ws.Range("N80").Value = -4996

# Row 83
$ws.Range("H83").Value = 3946.5386
$ws.Range("I83").Value = 4025.4167
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 20127.0835
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -15135.0835
$ws.Range("N83").Value = -24984

# Row 126
$ws.Range("H126").Value = 2450
$ws.Range("I126").Value = 3050
$ws.Range("J126").Value = 1850
$ws.Range("K126").Value = 9150
$ws.Range("L126").Value = 5550
$ws.Range("M126").Value = -6680
$ws.Range("N126").Value = -10490

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1291.2727
$ws.Range("I93").Value = 1166.6666
$ws.Range("J93").Value = 1852
$ws.Range("K93").Value = 1166.6666
$ws.Range("L93").Value = 1852
$ws.Range("M93").Value = 81.33339999999998
$ws.Range("N93").Value = -4348

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# Row 93
$ws.Range("H93").Value = 28254
$ws.Range("J93").Value = 28254
$ws.Range("L93").Value = 28254
$ws.Range("N93").Value = -33246

# Row 107
$ws.Range("H107").Value = 451.07693
$ws.Range("I107").Value = 457.1
$ws.Range("J107").Value = 431
$ws.Range("K107").Value = 1371.3
$ws.Range("L107").Value = 1293
$ws.Range("M107").Value = 548.6999999999998
$ws.Range("N107").Value = -5133

# Row 126
$ws.Range("H126").Value = 847.2727
$ws.Range("I126").Value = 565
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 1695
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = 775
$ws.Range("N126").Value = -9740
